# WSL (Wei Si Li) position statistics update: add a new date sheet "20191111"
# by copying the previous "20191110" sheet, updating its figures, and
# touching up selections on a couple of the other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet "20191110" (10th sheet) already has its C24 total formula in
#    this workbook; nothing to add there.
# ---------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item(10)

# ---------------------------------------------------------------------
# 2. Duplicate "20191110" to create the new "20191111" sheet right after it.
# ---------------------------------------------------------------------
$ws10.Copy($null, $ws10)
$ws11 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws11.Name = "20191111"

# ---------------------------------------------------------------------
# 3. Update the figures on the new "20191111" sheet with the day's data.
# ---------------------------------------------------------------------

# Distribution table (rows 1-14)
$ws11.Range("B1").Value = 979
$ws11.Range("C1").Value = 0.07
$ws11.Range("D1").Value = 0

$ws11.Range("B2").Value = 392
$ws11.Range("C2").Value = 0.02
$ws11.Range("D2").Value = 5

$ws11.Range("B3").Value = 454
$ws11.Range("C3").Value = 0.03
$ws11.Range("D3").Value = 15

$ws11.Range("B4").Value = 500
$ws11.Range("C4").Value = 0.03
$ws11.Range("D4").Value = 25

$ws11.Range("B5").Value = 497
$ws11.Range("C5").Value = 0.03
$ws11.Range("D5").Value = 35

$ws11.Range("B6").Value = 578
$ws11.Range("C6").Value = 0.04
$ws11.Range("D6").Value = 45

$ws11.Range("B7").Value = 773
$ws11.Range("C7").Value = 0.05
$ws11.Range("D7").Value = 55

$ws11.Range("B8").Value = 655
$ws11.Range("C8").Value = 0.04
$ws11.Range("D8").Value = 65

$ws11.Range("B9").Value = 810
$ws11.Range("C9").Value = 0.06
$ws11.Range("D9").Value = 75

$ws11.Range("B10").Value = 1004
$ws11.Range("C10").Value = 0.07
$ws11.Range("D10").Value = 85

$ws11.Range("B11").Value = 1563
$ws11.Range("C11").Value = 0.11
$ws11.Range("D11").Value = 95

$ws11.Range("B12").Value = 1619
$ws11.Range("C12").Value = 0.12
$ws11.Range("D12").Value = 100

$ws11.Range("B13").Value = 1113
$ws11.Range("C13").Value = 0.08
$ws11.Range("D13").Value = 100

$ws11.Range("B14").Value = 2171
$ws11.Range("C14").Value = 0.16

# Sentiment table (rows 20-23) - text labels differ slightly from the
# previous day's sheet ("看空" / "看平 (已选)" rather than "看空 (已选)" / "看平").
$ws11.Range("A20").Value = "看多"
$ws11.Range("B20").Value = 3809
$ws11.Range("C20").Value = 0.29

$ws11.Range("A21").Value = "看空"
$ws11.Range("B21").Value = 5541
$ws11.Range("C21").Value = 0.42

$ws11.Range("A22").Value = "看平 (已选)"
$ws11.Range("B22").Value = 1531
$ws11.Range("C22").Value = 0.11

$ws11.Range("A23").Value = "我是来给卫斯理打Call的"
$ws11.Range("B23").Value = 2173
$ws11.Range("C23").Value = 0.16

# ---------------------------------------------------------------------
# 4. Make the new sheet the active one and set its selection.
# ---------------------------------------------------------------------
$ws11.Activate()
$ws11.Range("E21").Select()

# ---------------------------------------------------------------------
# 5. Touch up selections left on the two sheets mentioned in the diff.
# ---------------------------------------------------------------------
$ws10.Activate()
$ws10.Range("H34").Select()

$ws9 = $wb.Worksheets.Item(9)
$ws9.Activate()
$ws9.Range("F16").Select()

# ---------------------------------------------------------------------
# 6. Restore the new sheet as the active tab (matches activeTab="10").
# ---------------------------------------------------------------------
$ws11.Activate()
